$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "Save" in H1, formatted like the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data column values (Save indicator 0/1)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1

$excel.CutCopyMode = 0
